# Auto-generated update of refreshed Universalis market-price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) across several
# Leve-profit worksheets, per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

# =========================== Sheet: ALC ===========================
$ws = $wb.Worksheets.Item("ALC")

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 459.30768
$ws.Range("I33").Value = 107.57895
$ws.Range("K33").Value = 107.57895
$ws.Range("M33").Value = 121.42105

# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 671.4
$ws.Range("I38").Value = 122.5
$ws.Range("J38").Value = 1037.3334
$ws.Range("K38").Value = 367.5
$ws.Range("L38").Value = 3112.0002
$ws.Range("M38").Value = 4.5
$ws.Range("N38").Value = -3856.0002

# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 164.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 4086.4092
$ws.Range("I74").Value = 4200.0713
$ws.Range("J74").Value = 3887.5
$ws.Range("K74").Value = 4200.0713
$ws.Range("L74").Value = 3887.5
$ws.Range("M74").Value = -3264.0713
$ws.Range("N74").Value = -5759.5

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 4086.4092
$ws.Range("I77").Value = 4200.0713
$ws.Range("J77").Value = 3887.5
$ws.Range("K77").Value = 21000.3565
$ws.Range("L77").Value = 19437.5
$ws.Range("M77").Value = -16320.3565
$ws.Range("N77").Value = -28797.5

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 28880.568
$ws.Range("I135").Value = 38334.703
$ws.Range("J135").Value = 3354.4
$ws.Range("K135").Value = 345012.327
$ws.Range("L135").Value = 30189.6
$ws.Range("M135").Value = -342477.327
$ws.Range("N135").Value = -35259.60000000001

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1667957.6
$ws.Range("I137").Value = 2326584
$ws.Range("J137").Value = 2020.6471
$ws.Range("K137").Value = 6979752
$ws.Range("L137").Value = 6061.9413
$ws.Range("M137").Value = -6977202
$ws.Range("N137").Value = -11161.9413

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1323.66
$ws.Range("I138").Value = 459.36
$ws.Range("J138").Value = 2187.96
$ws.Range("K138").Value = 1378.08
$ws.Range("L138").Value = 6563.88
$ws.Range("M138").Value = 3761.92
$ws.Range("N138").Value = -16843.88

# =========================== Sheet: ARM ===========================
$ws = $wb.Worksheets.Item("ARM")

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1268.55
$ws.Range("I32").Value = 887.51807
$ws.Range("J32").Value = 3128.8823
$ws.Range("K32").Value = 887.51807
$ws.Range("L32").Value = 3128.8823
$ws.Range("M32").Value = -600.51807
$ws.Range("N32").Value = -3702.8823

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 7413320.5
$ws.Range("I74").Value = 10041308
$ws.Range("K74").Value = 10041308
$ws.Range("M74").Value = -10040434

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 7413320.5
$ws.Range("I77").Value = 10041308
$ws.Range("K77").Value = 50206540
$ws.Range("M77").Value = -50202172

# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 34975
$ws.Range("J80").Value = 34966.668
$ws.Range("L80").Value = 34966.668
$ws.Range("N80").Value = -36962.668

# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 34975
$ws.Range("J83").Value = 34966.668
$ws.Range("L83").Value = 104900.004
$ws.Range("N83").Value = -114884.004

# =========================== Sheet: BSM ===========================
$ws = $wb.Worksheets.Item("BSM")

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1380.2858
$ws.Range("I107").Value = 1106.3077
$ws.Range("J107").Value = 2171.7778
$ws.Range("K107").Value = 1106.3077
$ws.Range("L107").Value = 2171.7778
$ws.Range("M107").Value = 813.6922999999999
$ws.Range("N107").Value = -6011.7778

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2276.8667
$ws.Range("I134").Value = 1335.5
$ws.Range("J134").Value = 4865.625
$ws.Range("K134").Value = 4006.5
$ws.Range("L134").Value = 14596.875
$ws.Range("M134").Value = -1471.5
$ws.Range("N134").Value = -19666.875

# Row 140: Ceremonial Teeth / Ra'Kaznar Twinfangs
$ws.Range("H140").Value = 53372
$ws.Range("J140").Value = 53372
$ws.Range("L140").Value = 53372
$ws.Range("N140").Value = -63732

# Row 141: Awl Dreams Come True / Ra'Kaznar Awl
$ws.Range("H141").Value = 46745.43
$ws.Range("J141").Value = 49160
$ws.Range("L141").Value = 49160
$ws.Range("N141").Value = -59520

# =========================== Sheet: CRP ===========================
$ws = $wb.Worksheets.Item("CRP")

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3202.361
$ws.Range("I31").Value = 1209.8276
$ws.Range("K31").Value = 1209.8276
$ws.Range("M31").Value = -914.8276000000001

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3202.361
$ws.Range("I34").Value = 1209.8276
$ws.Range("K34").Value = 1209.8276
$ws.Range("M34").Value = -1007.8276

# Row 51: Greenstone for Greenhorns / Jade Crook
$ws.Range("H51").Value = 24824.25
$ws.Range("J51").Value = 24824.25
$ws.Range("L51").Value = 24824.25
$ws.Range("N51").Value = -26296.25

# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 30845.25
$ws.Range("J59").Value = 30845.25
$ws.Range("L59").Value = 30845.25
$ws.Range("N59").Value = -33135.25

# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 9500

# Row 61: Incant Now, Think Later / Jade Crook
$ws.Range("H61").Value = 24824.25
$ws.Range("J61").Value = 24824.25
$ws.Range("L61").Value = 24824.25
$ws.Range("N61").Value = -25520.25

# Row 68: Do You Even String Bow / Holy Cedar Composite Bow
$ws.Range("H68").Value = 32196.666
$ws.Range("J68").Value = 35636
$ws.Range("L68").Value = 35636
$ws.Range("N68").Value = -37134

# Row 71: Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws.Range("H71").Value = 32196.666
$ws.Range("J71").Value = 35636
$ws.Range("L71").Value = 106908
$ws.Range("N71").Value = -114396

# =========================== Sheet: CUL ===========================
$ws = $wb.Worksheets.Item("CUL")

# Row 119: Super Dark Times / Risotto al Nero
$ws.Range("H119").Value = 2981

# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 15000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 45000
$ws.Range("N120").Value = -54676
$ws.Range("M120").ClearContents()

# =========================== Sheet: GSM ===========================
$ws = $wb.Worksheets.Item("GSM")

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3059.738
$ws.Range("J80").Value = 3370.3704
$ws.Range("L80").Value = 3370.3704
$ws.Range("N80").Value = -5366.3704

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3059.738
$ws.Range("J83").Value = 3370.3704
$ws.Range("L83").Value = 16851.852
$ws.Range("N83").Value = -26835.852

# Row 117: Birth Ring / Triplite Ring of Aiming
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 1523.2858
$ws.Range("I126").Value = 1150
$ws.Range("J126").Value = 1803.25
$ws.Range("K126").Value = 3450
$ws.Range("L126").Value = 5409.75
$ws.Range("M126").Value = -980
$ws.Range("N126").Value = -10349.75

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 66072.94
$ws.Range("I132").Value = 37223.605
$ws.Range("K132").Value = 111670.815
$ws.Range("M132").Value = -109140.815

# =========================== Sheet: WVR ===========================
$ws = $wb.Worksheets.Item("WVR")

# Row 41: Half Is the New Double / Linen Halfgloves
$ws.Range("H41").Value = 6758.4287
$ws.Range("I41").Value = 8121
$ws.Range("J41").Value = 6213.4
$ws.Range("K41").Value = 8121
$ws.Range("L41").Value = 6213.4
$ws.Range("M41").Value = -7731
$ws.Range("N41").Value = -6993.4

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1359.3077
$ws.Range("I126").Value = 1416
$ws.Range("J126").Value = 1047.5
$ws.Range("K126").Value = 4248
$ws.Range("L126").Value = 3142.5
$ws.Range("M126").Value = -1778
$ws.Range("N126").Value = -8082.5

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 51537.44
$ws.Range("I132").Value = 39146.23
$ws.Range("J132").Value = 123131.11
$ws.Range("K132").Value = 117438.69
$ws.Range("L132").Value = 369393.33
$ws.Range("M132").Value = -114908.69
$ws.Range("N132").Value = -374453.33

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 31656.83
$ws.Range("I136").Value = 18952.111
$ws.Range("J136").Value = 94025.45
$ws.Range("K136").Value = 56856.333
$ws.Range("L136").Value = 282076.35
$ws.Range("M136").Value = -54306.333
$ws.Range("N136").Value = -287176.35
